$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.599.96"
$ws.Range("E2").Value = "  +1.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.600.09"
$ws.Range("E3").Value = "  +1.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.39"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0911"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.828.84"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.599.71"
$ws.Range("E13").Value = "  +1.08%  "
$ws.Range("E14").Value = "  +3.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.605.02"
$ws.Range("E15").Value = "  +1.51%  "
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0693"
$ws.Range("E20").Value = "  +0.34%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.06"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  +0.43%  "
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.422.12"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  +1.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.88"
$ws.Range("E36").Value = "  +4.82%  "
$ws.Range("E38").Value = "  +0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0169"
$ws.Range("E39").Value = "  +3.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.545"
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.53"
$ws.Range("E41").Value = "  +6.93%  "
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("E43").Value = "  +5.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.808"
$ws.Range("E44").Value = "  +2.39%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.990"
$ws.Range("E46").Value = "  +16.85%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "66.32"
$ws.Range("E47").Value = "  +2.47%  "
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.739.14"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "86.22"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0104"
$ws.Range("E51").Value = "  +3.03%  "
